# Slide 5 ("ZoneTexte 4" shape) has a bullet paragraph whose sole run
# reads "diagramme de classes ". The author capitalised the first
# letter ("Diagramme de classes "), which causes PowerPoint to split
# the original run into three runs: "D" | "iagramme " | "de classes ".
# All three keep the same character formatting (accent3 / lumMod 50000).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.Trim() -eq "diagramme de classes") {
        $start = $para.Start

        # "d" -> "D"  (run 1, length 1)
        $run1 = $tr.Characters($start, 1)
        $run1.Text = "D"

        # re-assert "de classes " (run 3, length 11) to force PowerPoint
        # to split the remaining "iagramme de classes " text into
        # "iagramme " (run 2) + "de classes " (run 3)
        $run3 = $tr.Characters($start + 10, 11)
        $run3.Text = "de classes "
    }
}
